{"js": "// Update the EQ-5D-5L utility score table with the final re-run numbers.\n//\n// Table layout (0-indexed rows/cells):\n//   row 1: Characteristic | Before first symptoms\\nN = 168[1] | Today\\nN = 168[1] | p-value[2]\n//   row 2: EQ-5D-5L utility score | 76 (21) [73, 79] | 60 (23) [57, 64] | <0.001\n//   row 3: Missing | 23 | 22 | (blank)\n//\n// New values:\n//   N = 168  -> N = 159   (both \"Before first symptoms\" and \"Today\" header cells)\n//   76 (21) [73, 79]  -> 74 (21) [70, 77]\n//   60 (23) [57, 64]  -> 59 (22) [56, 63]\n//   23 (Missing, Before first symptoms) -> 1\n//   22 (Missing, Today)                 -> 0\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst headerRow = rows.items[1]; // \"Characteristic\" / N = 168 row\nconst meanRow = rows.items[2]; // \"EQ-5D-5L utility score\" row\nconst missingRow = rows.items[3]; // \"Missing\" row\n\nheaderRow.cells.load(\"items\");\nmeanRow.cells.load(\"items\");\nmissingRow.cells.load(\"items\");\nawait context.sync();\n\nasync function replaceInCell(cellBody, find, replace) {\n  const results = cellBody.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// \"N = 168\" -> \"N = 159\" in both the \"Before first symptoms\" and \"Today\" cells.\nawait replaceInCell(headerRow.cells.items[1].body, \"N = 168\", \"N = 159\");\nawait replaceInCell(headerRow.cells.items[2].body, \"N = 168\", \"N = 159\");\n\n// Mean (SD) [95% CI] values.\nawait replaceInCell(meanRow.cells.items[1].body, \"76 (21) [73, 79]\", \"74 (21) [70, 77]\");\nawait replaceInCell(meanRow.cells.items[2].body, \"60 (23) [57, 64]\", \"59 (22) [56, 63]\");\n\n// \"Missing\" counts.\nawait replaceInCell(missingRow.cells.items[1].body, \"23\", \"1\");\nawait replaceInCell(missingRow.cells.items[2].body, \"22\", \"0\");\n", "ps1": "# Update the EQ-5D-5L utility score table with the final re-run numbers.\n#\n# Table layout (1-indexed rows/cols, as COM exposes it):\n#   Row 2: Characteristic | Before first symptoms / N = 168[1] | Today / N = 168[1] | p-value[2]\n#   Row 3: EQ-5D-5L utility score | 76 (21) [73, 79] | 60 (23) [57, 64] | <0.001\n#   Row 4: Missing | 23 | 22 | (blank)\n#\n# New values:\n#   N = 168 -> N = 159   (both \"Before first symptoms\" and \"Today\" cells)\n#   76 (21) [73, 79] -> 74 (21) [70, 77]\n#   60 (23) [57, 64] -> 59 (22) [56, 63]\n#   23 (Missing, Before first symptoms) -> 1\n#   22 (Missing, Today)                 -> 0\n#\n# Note: several of these cells contain more than one run (e.g. a line break\n# plus a superscript footnote marker), so we can't just overwrite the whole\n# cell's Range.Text. Instead we locate the exact substring inside the cell\n# text and build a precise sub-Range (by absolute character offsets) so only\n# that substring is replaced, leaving surrounding runs/formatting untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfunction Replace-TextInCell($row, $col, $findText, $replaceText) {\n    $cell = $t.Cell($row, $col)\n    $cellRange = $cell.Range\n    $cellText = $cellRange.Text\n    $idx = $cellText.IndexOf($findText)\n    if ($idx -lt 0) {\n        throw \"Text '$findText' not found in cell ($row,$col)\"\n    }\n    $absStart = $cellRange.Start + $idx\n    $absEnd = $absStart + $findText.Length\n    $subRange = $d.Range($absStart, $absEnd)\n    $subRange.Text = $replaceText\n}\n\nReplace-TextInCell 2 2 \"N = 168\" \"N = 159\"\nReplace-TextInCell 2 3 \"N = 168\" \"N = 159\"\n\nReplace-TextInCell 3 2 \"76 (21) [73, 79]\" \"74 (21) [70, 77]\"\nReplace-TextInCell 3 3 \"60 (23) [57, 64]\" \"59 (22) [56, 63]\"\n\nReplace-TextInCell 4 2 \"23\" \"1\"\nReplace-TextInCell 4 3 \"22\" \"0\"\n"}
